{"js": "// Update the \"Sprint No.\" value from \"1\" to \"2\" and the \"Review Date\"\n// value from \"02/09/18\" to \"02/21/18\" in the checklist header table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// The header table's logical rows (accounting for merged/gridSpan cells)\n// are:\n//   Row 0: Project Name    | <value> | Project ID  | <value>\n//   Row 1: Reviewer's Name | <value> | Sprint No.  | <value>\n//   Row 2: Review Date     | <value (spans remaining columns)>\n//   Row 3: File Name (Source Code) | <value>\nconst sprintCell = table.getCellOrNullObject(1, 3);\nsprintCell.load(\"value\");\nawait context.sync();\n\nif (!sprintCell.isNullObject) {\n  const sprintResults = sprintCell.body.search(\"1\", { matchWholeWord: true });\n  sprintResults.load(\"items\");\n  await context.sync();\n  for (const r of sprintResults.items) {\n    r.insertText(\"2\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst reviewDateCell = table.getCellOrNullObject(2, 1);\nreviewDateCell.load(\"value\");\nawait context.sync();\n\nif (!reviewDateCell.isNullObject) {\n  const dateResults = reviewDateCell.body.search(\"02/09/18\", { matchWholeWord: false });\n  dateResults.load(\"items\");\n  await context.sync();\n  for (const r of dateResults.items) {\n    r.insertText(\"02/21/18\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the \"Sprint No.\" value from \"1\" to \"2\" and the \"Review Date\"\n# value from \"02/09/18\" to \"02/21/18\" in the checklist header table.\n\n$d = $word.ActiveDocument\n\n$table = $d.Tables.Item(1)\n\n# Logical rows of the header table (accounting for merged/gridSpan cells):\n#   Row 1: Project Name    | <value> | Project ID  | <value>\n#   Row 2: Reviewer's Name | <value> | Sprint No.  | <value>\n#   Row 3: Review Date     | <value (spans remaining columns)>\n#   Row 4: File Name (Source Code) | <value>\n\n$sprintCell = $table.Cell(2, 4)\n$sprintRange = $sprintCell.Range\n$find1 = $sprintRange.Find\n$find1.ClearFormatting()\n$find1.Text = \"1\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"2\"\n$find1.MatchWholeWord = $true\n$find1.Execute(\n    $find1.Text,\n    $false,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    $find1.Replacement.Text,\n    1\n)\n\n$dateCell = $table.Cell(3, 2)\n$dateRange = $dateCell.Range\n$find2 = $dateRange.Find\n$find2.ClearFormatting()\n$find2.Text = \"02/09/18\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"02/21/18\"\n$find2.MatchWholeWord = $false\n$find2.Execute(\n    $find2.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    $find2.Replacement.Text,\n    1\n)\n"}
